# Update "想去人数" (F column) figures on both the "展览" and "全部类型"
# worksheets, which mirror each other in this workbook.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F (想去人数)
$updates = @{
    3  = 1411
    4  = 167
    11 = 4709
    12 = 6969
    18 = 4172
    19 = 1003
    20 = 78
    21 = 70
    22 = 2748
    23 = 575
    25 = 176
    28 = 405
    32 = 1050
    33 = 72
    34 = 619
    35 = 89
    36 = 553
    37 = 7
    43 = 22
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
